# Update the "Pais" (countries) worksheet with refreshed COVID-19 figures.
# The source feed was re-pulled; the table is sorted descending by "Casos
# totales" so a handful of neighbouring countries swap rank (their row
# position / country name) while others simply get new totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Header: last-updated timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Septiembre de 2020 a las 05:05"

# --- Row 8: Peru (figures refreshed, no rank change) ---
$ws.Range("A8").Value = "Peru"
$ws.Range("B8").Value = 750098
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 594513
$ws.Range("E8").Value = 124439
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 31146

# --- Row 30: Bolivia (figures refreshed, no rank change) ---
$ws.Range("A30").Value = "Bolivia"
$ws.Range("B30").Value = 129419
$ws.Range("C30").Value = 547
$ws.Range("D30").Value = 87716
$ws.Range("E30").Value = 34192
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 33
$ws.Range("H30").Value = 7511

# --- Row 34: Kazajistan (figures refreshed, no rank change) ---
$ws.Range("A34").Value = "Kazajistan"
$ws.Range("B34").Value = 107134
$ws.Range("C34").Value = 78
$ws.Range("D34").Value = 101610
$ws.Range("E34").Value = 3853
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 1671

# --- Rows 38-39: Belgica overtakes Kuwait ---
$ws.Range("A38").Value = "Belgica"
$ws.Range("B38").Value = 97976
$ws.Range("C38").Value = 2028
$ws.Range("D38").Value = 18854
$ws.Range("E38").Value = 69186
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 9936

$ws.Range("A39").Value = "Kuwait"
$ws.Range("B39").Value = 97824
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 87911
$ws.Range("E39").Value = 9338
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 575

# --- Row 50: Honduras (figures refreshed, no rank change) ---
$ws.Range("A50").Value = "Honduras"
$ws.Range("B50").Value = 70120
$ws.Range("C50").Value = 460
$ws.Range("D50").Value = 20677
$ws.Range("E50").Value = 47321
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 20
$ws.Range("H50").Value = 2122

# --- Row 53: Venezuela (figures refreshed, no rank change) ---
$ws.Range("A53").Value = "Venezuela"
$ws.Range("B53").Value = 64284
$ws.Range("C53").Value = 0
$ws.Range("D53").Value = 52564
$ws.Range("E53").Value = 11200
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 520

# --- Row 76: Australia (figures refreshed, no rank change) ---
$ws.Range("A76").Value = "Australia"
$ws.Range("B76").Value = 26861
$ws.Range("C76").Value = 48
$ws.Range("D76").Value = 23793
$ws.Range("E76").Value = 2231
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 5
$ws.Range("H76").Value = 837

# --- Rows 126-128: Birmania jumps ahead of Jordania and Eslovenia ---
$ws.Range("A126").Value = "Birmania"
$ws.Range("B126").Value = 4299
$ws.Range("C126").Value = 256
$ws.Range("D126").Value = 944
$ws.Range("E126").Value = 3294
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 61

$ws.Range("A127").Value = "Jordania"
$ws.Range("B127").Value = 4131
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 2415
$ws.Range("E127").Value = 1690
$ws.Range("F127").Value = 0
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 26

$ws.Range("A128").Value = "Eslovenia"
$ws.Range("B128").Value = 4058
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 2897
$ws.Range("E128").Value = 1025
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 136

# --- Rows 158-159: Belice overtakes Republica de Chipre ---
$ws.Range("A158").Value = "Belice"
$ws.Range("B158").Value = 1567
$ws.Range("C158").Value = 31
$ws.Range("D158").Value = 742
$ws.Range("E158").Value = 806
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 19

$ws.Range("A159").Value = "Republica de Chipre"
$ws.Range("B159").Value = 1558
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 1282
$ws.Range("E159").Value = 254
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 22

# --- Row 173: San Martin (Parte Holandesa) (figures refreshed, no rank change) ---
$ws.Range("A173").Value = "San Martin (Parte Holandesa)"
$ws.Range("B173").Value = 565
$ws.Range("C173").Value = 8
$ws.Range("D173").Value = 477
$ws.Range("E173").Value = 69
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 19

# --- Rows 214-215: Islas Malvinas overtakes Montserrat (tied totals) ---
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1
